$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (strikeouts) values replacing the old "Strike#" values in column G
$values = @{
    2  = 2
    3  = 1
    4  = 2
    5  = 1
    6  = 1
    7  = 1
    8  = 1
    9  = 6
    10 = 1
    11 = 7
    12 = 3
    13 = 4
    14 = 7
    15 = 2
    16 = 4
    17 = 1
    18 = 6
    19 = 3
    20 = 4
    21 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
